# Auto-generated script to apply scheduled-runner profit recalculation updates
# to the Belias_Profits workbook, sheet by sheet, cell by cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 15887.5
$ws.Range("J7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("N7").Value = -22224
$ws.Range("H14").Value = 15887.5
$ws.Range("J14").Value = 22000
$ws.Range("L14").Value = 22000
$ws.Range("N14").Value = -22382
$ws.Range("H46").Value = 19297
$ws.Range("I46").Value = 3195.5
$ws.Range("J46").Value = 51500
$ws.Range("K46").Value = 9586.5
$ws.Range("L46").Value = 154500
$ws.Range("M46").Value = -9467.5
$ws.Range("N46").Value = -154738
$ws.Range("H51").Value = 3490.1365
$ws.Range("I51").Value = 1633.5
$ws.Range("J51").Value = 4186.375
$ws.Range("K51").Value = 1633.5
$ws.Range("L51").Value = 4186.375
$ws.Range("M51").Value = -1149.5
$ws.Range("N51").Value = -5154.375
$ws.Range("H60").Value = 19297
$ws.Range("I60").Value = 3195.5
$ws.Range("J60").Value = 51500
$ws.Range("K60").Value = 9586.5
$ws.Range("L60").Value = 154500
$ws.Range("M60").Value = -9102.5
$ws.Range("N60").Value = -155468
$ws.Range("H64").Value = 3942.8572
$ws.Range("I64").Value = 3900
$ws.Range("J64").Value = 3950
$ws.Range("K64").Value = 3900
$ws.Range("L64").Value = 3950
$ws.Range("M64").Value = -3652
$ws.Range("N64").Value = -4446
$ws.Range("H67").Value = 3942.8572
$ws.Range("I67").Value = 3900
$ws.Range("J67").Value = 3950
$ws.Range("K67").Value = 3900
$ws.Range("L67").Value = 3950
$ws.Range("M67").Value = -3042
$ws.Range("N67").Value = -5666
$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H93").Value = 29950
$ws.Range("J93").Value = 29950
$ws.Range("L93").Value = 29950
$ws.Range("N93").Value = -34942
$ws.Range("H98").Value = 1256.2646
$ws.Range("I98").Value = 725.2857
$ws.Range("K98").Value = 725.2857
$ws.Range("M98").Value = 772.7143
$ws.Range("H100").Value = 2061.739
$ws.Range("I100").Value = 1306.25
$ws.Range("J100").Value = 2464.6667
$ws.Range("K100").Value = 1306.25
$ws.Range("L100").Value = 2464.6667
$ws.Range("M100").Value = -765.25
$ws.Range("N100").Value = -3546.6667
$ws.Range("H116").Value = 4819.3125
$ws.Range("I116").Value = 1643.1428
$ws.Range("J116").Value = 7289.6665
$ws.Range("K116").Value = 1643.1428
$ws.Range("L116").Value = 7289.6665
$ws.Range("M116").Value = 1798.8572
$ws.Range("N116").Value = -14173.6665
$ws.Range("H122").Value = 1256.2646
$ws.Range("I122").Value = 725.2857
$ws.Range("K122").Value = 2175.8571
$ws.Range("M122").Value = 274.1428999999998
$ws.Range("H127").Value = 125005010
$ws.Range("I127").Value = 200000480
$ws.Range("K127").Value = 600001440
$ws.Range("M127").Value = -599996480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2369
$ws.Range("I2").Value = 2168.2222
$ws.Range("J2").Value = 2627.1428
$ws.Range("K2").Value = 2168.2222
$ws.Range("L2").Value = 2627.1428
$ws.Range("M2").Value = -2055.2222
$ws.Range("N2").Value = -2853.1428
$ws.Range("H10").Value = 18876.5
$ws.Range("I10").Value = 1833.6666
$ws.Range("J10").Value = 70005
$ws.Range("K10").Value = 1833.6666
$ws.Range("L10").Value = 70005
$ws.Range("M10").Value = -1663.6666
$ws.Range("N10").Value = -70345
$ws.Range("H45").Value = 1505.9395
$ws.Range("I45").Value = 1181.2667
$ws.Range("J45").Value = 1776.5
$ws.Range("K45").Value = 1181.2667
$ws.Range("L45").Value = 1776.5
$ws.Range("M45").Value = -804.2666999999999
$ws.Range("N45").Value = -2530.5
$ws.Range("H88").Value = 2371.4285
$ws.Range("J88").Value = 2200
$ws.Range("L88").Value = 2200
$ws.Range("N88").Value = -3012
$ws.Range("H91").Value = 2371.4285
$ws.Range("J91").Value = 2200
$ws.Range("L91").Value = 2200
$ws.Range("N91").Value = -5008
$ws.Range("H116").Value = 2369
$ws.Range("I116").Value = 2168.2222
$ws.Range("J116").Value = 2627.1428
$ws.Range("K116").Value = 2168.2222
$ws.Range("L116").Value = 2627.1428
$ws.Range("M116").Value = 125.7777999999998
$ws.Range("N116").Value = -7215.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2369
$ws.Range("I3").Value = 2168.2222
$ws.Range("J3").Value = 2627.1428
$ws.Range("K3").Value = 2168.2222
$ws.Range("L3").Value = 2627.1428
$ws.Range("M3").Value = -2054.2222
$ws.Range("N3").Value = -2855.1428
$ws.Range("H20").Value = 1585.9524
$ws.Range("I20").Value = 1372.7693
$ws.Range("J20").Value = 1932.375
$ws.Range("K20").Value = 1372.7693
$ws.Range("L20").Value = 1932.375
$ws.Range("M20").Value = -1125.7693
$ws.Range("N20").Value = -2426.375
$ws.Range("H107").Value = 1625.4642
$ws.Range("I107").Value = 1535.7142
$ws.Range("J107").Value = 1894.7142
$ws.Range("K107").Value = 1535.7142
$ws.Range("L107").Value = 1894.7142
$ws.Range("M107").Value = 384.2858000000001
$ws.Range("N107").Value = -5734.7142
$ws.Range("H141").Value = 99000
$ws.Range("J141").Value = 99000
$ws.Range("L141").Value = 99000
$ws.Range("N141").Value = -109360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 931
$ws.Range("I16").Value = 762
$ws.Range("J16").Value = 1032.4
$ws.Range("K16").Value = 762
$ws.Range("L16").Value = 1032.4
$ws.Range("M16").Value = -475
$ws.Range("N16").Value = -1606.4
$ws.Range("H31").Value = 27029158
$ws.Range("I31").Value = 71429540
$ws.Range("J31").Value = 2840.9565
$ws.Range("K31").Value = 71429540
$ws.Range("L31").Value = 2840.9565
$ws.Range("M31").Value = -71429245
$ws.Range("N31").Value = -3430.9565
$ws.Range("H34").Value = 27029158
$ws.Range("I34").Value = 71429540
$ws.Range("J34").Value = 2840.9565
$ws.Range("K34").Value = 71429540
$ws.Range("L34").Value = 2840.9565
$ws.Range("M34").Value = -71429338
$ws.Range("N34").Value = -3244.9565
$ws.Range("H113").Value = 931
$ws.Range("I113").Value = 762
$ws.Range("J113").Value = 1032.4
$ws.Range("K113").Value = 762
$ws.Range("L113").Value = 1032.4
$ws.Range("M113").Value = 1408
$ws.Range("N113").Value = -5372.4
$ws.Range("H134").Value = 55002960
$ws.Range("I134").Value = 5558566
$ws.Range("J134").Value = 500002500
$ws.Range("K134").Value = 16675698
$ws.Range("L134").Value = 1500007500
$ws.Range("M134").Value = -16673163
$ws.Range("N134").Value = -1500012570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5875
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 6500
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 19500
$ws.Range("M3").Value = -11888
$ws.Range("N3").Value = -19724
$ws.Range("H56").Value = 5079.1304
$ws.Range("I56").Value = 5079.1304
$ws.Range("K56").Value = 5079.1304
$ws.Range("M56").Value = -4549.1304
$ws.Range("H68").Value = 871.7143
$ws.Range("I68").Value = 588
$ws.Range("K68").Value = 1764
$ws.Range("M68").Value = -953
$ws.Range("H71").Value = 871.7143
$ws.Range("I71").Value = 588
$ws.Range("K71").Value = 5292
$ws.Range("M71").Value = -1236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2225.875
$ws.Range("I9").Value = 972.4286
$ws.Range("J9").Value = 11000
$ws.Range("K9").Value = 972.4286
$ws.Range("L9").Value = 11000
$ws.Range("M9").Value = -802.4286
$ws.Range("N9").Value = -11340
$ws.Range("H70").Value = 4968.86
$ws.Range("I70").Value = 5085.6665
$ws.Range("J70").Value = 4884.276
$ws.Range("K70").Value = 5085.6665
$ws.Range("L70").Value = 4884.276
$ws.Range("M70").Value = -4815.6665
$ws.Range("N70").Value = -5424.276
$ws.Range("H73").Value = 4968.86
$ws.Range("I73").Value = 5085.6665
$ws.Range("J73").Value = 4884.276
$ws.Range("K73").Value = 5085.6665
$ws.Range("L73").Value = 4884.276
$ws.Range("M73").Value = -4149.6665
$ws.Range("N73").Value = -6756.276
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H102").Value = 1435.1111
$ws.Range("I102").Value = 1155.6923
$ws.Range("J102").Value = 2161.6
$ws.Range("K102").Value = 1155.6923
$ws.Range("L102").Value = 2161.6
$ws.Range("M102").Value = 466.3077000000001
$ws.Range("N102").Value = -5405.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1612
$ws.Range("I100").Value = 1714
$ws.Range("K100").Value = 3428
$ws.Range("M100").Value = -2887
$ws.Range("H132").Value = 3484.5312
$ws.Range("I132").Value = 3929.1428
$ws.Range("J132").Value = 2635.7273
$ws.Range("K132").Value = 11787.4284
$ws.Range("L132").Value = 7907.1819
$ws.Range("M132").Value = -9257.428400000001
$ws.Range("N132").Value = -12967.1819
